# Update row 4 (patient "Aristoteles") values on the active worksheet
# to reflect corrected calculations (commit: "changed is operator to == (dumbass)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 109
$ws.Range("H4").Value = 1144
$ws.Range("I4").Value = 1015
$ws.Range("J4").Value = 109
$ws.Range("K4").Value = 223
$ws.Range("L4").Value = 271
$ws.Range("M4").Value = 581
$ws.Range("N4").Value = 652
$ws.Range("O4").Value = 1015
$ws.Range("P4").Value = 1144
$ws.Range("Q4").Value = 1246
$ws.Range("R4").Value = 1294
$ws.Range("V4").Value = -17.88
$ws.Range("W4").Value = 32.25
